# Thesis bibliography cleanup:
#  - remove four duplicate/redundant reference rows
#  - re-sort the reference list A -> Z

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

# Distinctive prefixes identifying the four rows to remove. Each of these
# references appears twice in the sheet; we keep the first occurrence
# (lowest row number) and delete every later duplicate.
$targets = @(
    "Garrison, E., & Marth, G. (2012). Haplotype-based variant detection from short-read sequencing.",
    "Sutskever, I., Martens, J., Dahl, G. E., & Hinton, G. E. (2013). On the importance of initialization and momentum in deep learning.",
    "Chawla, N. V. (2005). Data mining for imbalanced datasets: An overview.",
    "Zhang, H., Gao, J., Zhao, Z., Li, M., & Liu, C. (2014). Clinical implications of SPRR1A"
)

$rowsToDelete = New-Object System.Collections.ArrayList

foreach ($prefix in $targets) {
    $matches = New-Object System.Collections.ArrayList
    for ($r = 2; $r -le $lastRow; $r++) {
        $val = $ws.Cells.Item($r, 1).Value2
        if ($val -eq $null) { continue }
        if ($val.ToString().StartsWith($prefix)) {
            [void]$matches.Add($r)
        }
    }
    # Keep the first occurrence, delete the rest.
    for ($i = 1; $i -lt $matches.Count; $i++) {
        [void]$rowsToDelete.Add($matches[$i])
    }
}

# Delete from the bottom up so earlier row numbers stay valid.
$sortedDesc = $rowsToDelete | Sort-Object -Descending -Unique
foreach ($r in $sortedDesc) {
    $ws.Rows.Item($r).Delete()
}

# Re-sort the remaining reference list (column A, header in row 1).
$used = $ws.UsedRange
$lastRow = $used.Rows.Count
$sortRange = $ws.Range("A2:A" + $lastRow)

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($sortRange)
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlNo
$ws.Sort.Apply()
